$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab title to reflect new "through" date
$ws.Name = "Through 2022-02-21"

# Update the header label cell (I1) which uses the shared string "2022 (through 02-20)"
$ws.Range("I1").Value = "2022 (through 02-21)"

# Update February 2022 value (row 3 = February)
$ws.Range("I3").Value = 106

# Update Total 2022 value (row 14 = Total)
$ws.Range("I14").Value = 265
